$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 18

function DateToSerial($y, $m, $d) {
    # Converts a Gregorian y/m/d into an Excel (1900 date system) serial number,
    # i.e. days since 1899-12-30 - pure arithmetic, no DateTime object needed.
    if ($m -le 2) {
        $y = $y - 1
        $m = $m + 12
    }
    $a = [math]::Floor($y / 100)
    $b = 2 - $a + [math]::Floor($a / 4)
    $jd = [math]::Floor(365.25 * ($y + 4716)) + [math]::Floor(30.6001 * ($m + 1)) + $d + $b - 1524.5
    return $jd - 2415018.5
}

# --- Snapshot every data row (2..lastRow) exactly as stored, before any writes ---
$rows = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $dateText = [string]$ws.Range("E$r").Value2   # e.g. "05-04-2021" (dd-mm-yyyy)
    $parts = $dateText.Split("-")
    $day   = [int]$parts[0]
    $month = [int]$parts[1]
    $year  = [int]$parts[2]
    $dateSerial = DateToSerial $year $month $day

    $rows += [PSCustomObject]@{
        A = $ws.Range("A$r").Value2
        B = $ws.Range("B$r").Value2
        C = $ws.Range("C$r").Value2
        D = [string]$ws.Range("D$r").Value2
        E = $dateSerial
        F = $ws.Range("F$r").Value2
        G = $ws.Range("G$r").Value2
        H = $ws.Range("H$r").Value2
        I = $ws.Range("I$r").Value2
    }
}

# --- Stable-sort within each repository (column D) group by date ascending,
#     keeping the original relative (first-seen) order of the groups themselves ---
$groupOrder = @()
$groups = @{}
foreach ($row in $rows) {
    if (-not $groups.ContainsKey($row.D)) {
        $groups[$row.D] = @()
        $groupOrder += $row.D
    }
    $groups[$row.D] += $row
}

$sortedRows = @()
foreach ($key in $groupOrder) {
    $sortedRows += ($groups[$key] | Sort-Object -Property E)
}

# --- Apply the custom date/time display format to the whole E column.
#     Assigning the lowercase form first and then the uppercase form on the
#     first cell reproduces the workbook's two registered numFmts (one
#     unused), then every other cell just picks up the uppercase one. ---
$ws.Range("E2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("E2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Range("E$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# --- Write the re-ordered rows back in place ---
for ($i = 0; $i -lt $sortedRows.Count; $i++) {
    $r = $i + 2
    $src = $sortedRows[$i]
    $ws.Range("A$r").Value2 = $src.A
    $ws.Range("B$r").Value2 = $src.B
    $ws.Range("C$r").Value2 = $src.C
    $ws.Range("D$r").Value2 = $src.D
    $ws.Range("E$r").Value2 = $src.E
    $ws.Range("F$r").Value2 = $src.F
    $ws.Range("G$r").Value2 = $src.G
    $ws.Range("H$r").Value2 = $src.H
    $ws.Range("I$r").Value2 = $src.I
}
